# 自动更新Excel文件 - applies the day's roster refresh:
#   * rows 2-79: updated "剩余" (E) / "开始时间" (F) counters as the cycle
#     advanced, with a few rows also getting a corrected "总天" (D) and/or
#     a cleared "备注1" (G) note;
#   * rows 80-93 (newly added stores) appended at the bottom;
#   * the sheet's used-range <dimension> grows from A1:I79 to A1:I93
#     automatically because Excel recomputes it from the populated cells.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Per-row field updates for the existing 78 data rows (sheet rows 2-79).
#    Only the columns that actually changed are listed for each row:
#      D = 总天 (total days), E = 剩余 (remaining), F = 开始时间 (start date),
#      G = 备注1 (note 1, cleared to blank on a few rows).
# ---------------------------------------------------------------------------
$rowUpdates = @(
    @{ Row=2;  E=7;  F=20250926 },
    @{ Row=3;  E=7;  F=20250926 },
    @{ Row=4;  E=7;  F=20250926 },
    @{ Row=5;  E=7;  F=20250926 },
    @{ Row=6;  D=10; E=7;  F=20250926; G=$null },
    @{ Row=7;  E=7;  F=20250926 },
    @{ Row=8;  E=7;  F=20250926 },
    @{ Row=9;  E=7;  F=20250926 },
    @{ Row=10; E=4;  F=20250926; G=$null },
    @{ Row=11; E=7;  F=20250926 },
    @{ Row=12; D=10; E=7;  F=20250926; G=$null },
    @{ Row=13; E=7;  F=20250926 },
    @{ Row=14; E=7;  F=20250926 },
    @{ Row=15; E=7;  F=20250926 },
    @{ Row=16; E=9 },
    @{ Row=17; E=6;  F=20250925 },
    @{ Row=18; E=10; F=20250929 },
    @{ Row=19; E=10; F=20250929 },
    @{ Row=20; E=10; F=20250929 },
    @{ Row=21; E=10; F=20250929 },
    @{ Row=22; E=10; F=20250929 },
    @{ Row=23; E=6;  F=20250925 },
    @{ Row=24; E=6;  F=20250925 },
    @{ Row=25; E=2;  F=20250921 },
    @{ Row=26; E=6;  F=20250925 },
    @{ Row=27; E=1 },
    @{ Row=28; E=10; F=20250929 },
    @{ Row=29; E=10; F=20250929 },
    @{ Row=30; E=10; F=20250929 },
    @{ Row=31; E=10; F=20250929 },
    @{ Row=32; E=10; F=20250929 },
    @{ Row=33; E=10; F=20250929 },
    @{ Row=34; E=10; F=20250929 },
    @{ Row=35; E=10; F=20250929 },
    @{ Row=36; E=2 },
    @{ Row=37; E=2 },
    @{ Row=38; E=2 },
    @{ Row=39; E=2 },
    @{ Row=40; F=20250929 },
    @{ Row=41; F=20250929 },
    @{ Row=42; E=2 },
    @{ Row=43; D=10; E=2;  F=20250921 },
    @{ Row=44; F=20250929 },
    @{ Row=45; E=6;  F=20250925 },
    @{ Row=46; D=7;  E=7;  F=20250929 },
    @{ Row=47; E=4 },
    @{ Row=48; F=20250929 },
    @{ Row=49; E=7;  F=20250929 },
    @{ Row=50; E=5;  F=20250924 },
    @{ Row=51; E=5;  F=20250924 },
    @{ Row=52; E=5;  F=20250924 },
    @{ Row=53; E=5;  F=20250924 },
    @{ Row=54; E=5;  F=20250924 },
    @{ Row=55; E=5;  F=20250924 },
    @{ Row=56; D=10; E=5;  F=20250924 },
    @{ Row=57; E=5;  F=20250924 },
    @{ Row=58; D=10; E=9;  F=20250928 },
    @{ Row=59; E=9 },
    @{ Row=60; E=9 },
    @{ Row=61; E=7;  F=20250929 },
    @{ Row=62; E=9 },
    @{ Row=63; E=9 },
    @{ Row=64; E=9 },
    @{ Row=65; E=9 },
    @{ Row=66; D=10; E=10; F=20250929 },
    @{ Row=67; D=10; E=10; F=20250929 },
    @{ Row=68; F=20250929 },
    @{ Row=69; F=20250929 },
    @{ Row=70; E=1 },
    @{ Row=71; E=1 },
    @{ Row=72; E=1 },
    @{ Row=73; E=1 },
    @{ Row=74; E=1 },
    @{ Row=75; E=1 },
    @{ Row=76; E=1 },
    @{ Row=77; E=2 },
    @{ Row=78; E=2 }
)

foreach ($update in $rowUpdates) {
    $r = $update.Row
    if ($update.ContainsKey('D')) { $ws.Cells.Item($r, 4).Value = $update.D }
    if ($update.ContainsKey('E')) { $ws.Cells.Item($r, 5).Value = $update.E }
    if ($update.ContainsKey('F')) { $ws.Cells.Item($r, 6).Value = $update.F }
    if ($update.ContainsKey('G')) { $ws.Cells.Item($r, 7).Value = $update.G }
}

# ---------------------------------------------------------------------------
# 2. Newly added stores, appended as sheet rows 80-93
#    (行号 79-92): 行号, 店铺名称, 地址, 总天, 剩余, 开始时间, 备注1, 备注2, 备注3
# ---------------------------------------------------------------------------
$newRows = @(
    @(79, '周广平特色',           '大湖大街',  10, 4, 20250923, '大桶1个', $null,     $null),
    @(80, '食惠坊',               '大湖大街',  10, 4, 20250923, $null,     '小桶1个', $null),
    @(81, '味满堂',               '大湖大街',  10, 4, 20250923, $null,     '小桶1个', $null),
    @(82, '社区家常菜',           '大湖大街',  10, 4, 20250923, $null,     '小桶1个', $null),
    @(83, '淮扬面馆',             '崇义路',    10, 4, 20250923, '大桶1个', $null,     $null),
    @(84, '金陵水饺城',           '平阳东路',  10, 4, 20250923, $null,     '小桶1个', $null),
    @(85, '大路全羊',             '崇义路',    10, 4, 20250923, '大桶1个', $null,     $null),
    @(86, '范大碗',               '平阳路',    7,  2, 20250924, '大桶4个', $null,     $null),
    @(87, '微渔坊',               '峄山路',    7,  5, 20250927, $null,     '小桶2个', $null),
    @(88, '澳龙特色餐厅',         '崇义路',    7,  5, 20250927, $null,     '小桶4个', $null),
    @(89, '福源酒家',             '峄山路',    7,  6, 20250928, $null,     '小桶2个', $null),
    @(90, '九道菜',               '峄山路',    7,  6, 20250928, '大桶1个', $null,     $null),
    @(91, '百大生态园',           '峄山路',    7,  6, 20250928, '大桶5个', $null,     $null),
    @(92, '顺合庭私厨大湖大街',   $null,       10, 4, 20250923, '大桶2个', $null,     $null)
)

$startRow = 80
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $values = $newRows[$i]
    for ($c = 1; $c -le $values.Count; $c++) {
        $ws.Cells.Item($r, $c).Value = $values[$c - 1]
    }
}
